$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.293731927871704
$ws.Range("B1").Value = 2.070857286453247
$ws.Range("C1").Value = 4.679467678070068
$ws.Range("D1").Value = 3.445167779922485
$ws.Range("E1").Value = 1.372546434402466
